# The commit updates the "DBD" sheet's data-dictionary rows for the
# CreateDate and LastUpdate fields: their SQL data type ("形態") column
# changes from "DATE" to "TIMESTAMP".
#
#   Row 14 -> CreateDate      (D14: DATE -> TIMESTAMP)
#   Row 16 -> LastUpdate      (D16: DATE -> TIMESTAMP)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("D14").Value = "TIMESTAMP"
$ws.Range("D16").Value = "TIMESTAMP"
